# Applies updated crypto price/volume data to the active worksheet.
# Mirrors the commit "Updated cryptos list ... with GitHub Actions":
# the Price (column D) and Volume(1h) (column E) columns are refreshed
# for the coin rows; the underlying data must stay plain text (as the
# original cells use inline/shared strings, e.g. "29.992.71", not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"        # force text so "29.992.71" is not parsed as a number
$cell.Value = "29.992.71"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E2").Value = "  -0.39%  "

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"        # force text so "1.868.53" is not parsed as a number
$cell.Value = "1.868.53"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E3").Value = "  -2.69%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"        # force text so "319.02" is not parsed as a number
$cell.Value = "319.02"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E5").Value = "  -2.44%  "

# Row 6
$ws.Range("E6").Value = "  -0.04%  "

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"        # force text so "0.5090" is not parsed as a number
$cell.Value = "0.5090"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E7").Value = "  -1.31%  "

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"        # force text so "0.3922" is not parsed as a number
$cell.Value = "0.3922"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E8").Value = "  -2.19%  "

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"        # force text so "0.08157" is not parsed as a number
$cell.Value = "0.08157"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E9").Value = "  -3.48%  "

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"        # force text so "42.03" is not parsed as a number
$cell.Value = "42.03"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E10").Value = "  -1.75%  "

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"        # force text so "1.088" is not parsed as a number
$cell.Value = "1.088"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E11").Value = "  -3.03%  "

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"        # force text so "22.72" is not parsed as a number
$cell.Value = "22.72"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E12").Value = "  +4.67%  "

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"        # force text so "1.862.73" is not parsed as a number
$cell.Value = "1.862.73"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E13").Value = "  -3.05%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"        # force text so "6.238" is not parsed as a number
$cell.Value = "6.238"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E14").Value = "  -1.75%  "

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"        # force text so "7.146" is not parsed as a number
$cell.Value = "7.146"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E15").Value = "  -2.82%  "

# Row 16
$ws.Range("E16").Value = "  -0.04%  "

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"        # force text so "91.54" is not parsed as a number
$cell.Value = "91.54"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E17").Value = "  -4.82%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"        # force text so "0.00001077" is not parsed as a number
$cell.Value = "0.00001077"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E18").Value = "  -3.53%  "

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"        # force text so "0.06362" is not parsed as a number
$cell.Value = "0.06362"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E19").Value = "  -5.62%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"        # force text so "17.84" is not parsed as a number
$cell.Value = "17.84"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E20").Value = "  -1.38%  "

# Row 21
$ws.Range("E21").Value = "  +0.10%  "

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"        # force text so "29.982.15" is not parsed as a number
$cell.Value = "29.982.15"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E22").Value = "  -0.49%  "

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"        # force text so "5.783" is not parsed as a number
$cell.Value = "5.783"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E23").Value = "  -4.59%  "

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"        # force text so "11.06" is not parsed as a number
$cell.Value = "11.06"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E24").Value = "  -1.44%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"        # force text so "2.187" is not parsed as a number
$cell.Value = "2.187"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E25").Value = "  -0.72%  "

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"        # force text so "2.081.25" is not parsed as a number
$cell.Value = "2.081.25"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E26").Value = "  -2.82%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"        # force text so "161.15" is not parsed as a number
$cell.Value = "161.15"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E27").Value = "  +0.29%  "

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"        # force text so "20.89" is not parsed as a number
$cell.Value = "20.89"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E28").Value = "  -0.88%  "

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"        # force text so "2.215" is not parsed as a number
$cell.Value = "2.215"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E29").Value = "  -9.88%  "

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"        # force text so "126.77" is not parsed as a number
$cell.Value = "126.77"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E30").Value = "  -1.98%  "

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"        # force text so "1.048" is not parsed as a number
$cell.Value = "1.048"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E31").Value = "  -2.56%  "

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"        # force text so "0.1031" is not parsed as a number
$cell.Value = "0.1031"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E32").Value = "  -2.73%  "

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"        # force text so "5.880" is not parsed as a number
$cell.Value = "5.880"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E33").Value = "  -3.25%  "

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"        # force text so "3.731" is not parsed as a number
$cell.Value = "3.731"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E34").Value = "  +1.82%  "

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"        # force text so "0.02408" is not parsed as a number
$cell.Value = "0.02408"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E35").Value = "  -4.23%  "

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"        # force text so "5.193" is not parsed as a number
$cell.Value = "5.193"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E36").Value = "  -0.26%  "

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"        # force text so "0.06307" is not parsed as a number
$cell.Value = "0.06307"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E37").Value = "  -4.45%  "

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"        # force text so "0.2134" is not parsed as a number
$cell.Value = "0.2134"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E38").Value = "  -3.93%  "

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"        # force text so "1.172" is not parsed as a number
$cell.Value = "1.172"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E39").Value = "  -5.19%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"        # force text so "8.463" is not parsed as a number
$cell.Value = "8.463"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E40").Value = "  -6.16%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"        # force text so "0.6271" is not parsed as a number
$cell.Value = "0.6271"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E41").Value = "  -4.22%  "

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"        # force text so "1.204" is not parsed as a number
$cell.Value = "1.204"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E42").Value = "  -3.08%  "

# Row 43
$ws.Range("E43").Value = "  -1.48%  "

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"        # force text so "0.9998" is not parsed as a number
$cell.Value = "0.9998"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"        # force text so "0.5873" is not parsed as a number
$cell.Value = "0.5873"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E45").Value = "  -4.31%  "

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"        # force text so "12.88" is not parsed as a number
$cell.Value = "12.88"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E46").Value = "  -2.45%  "

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"        # force text so "3.620" is not parsed as a number
$cell.Value = "3.620"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E47").Value = "  -3.85%  "

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"        # force text so "1.979" is not parsed as a number
$cell.Value = "1.979"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E48").Value = "  -3.67%  "

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"        # force text so "122.55" is not parsed as a number
$cell.Value = "122.55"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E49").Value = "  -2.55%  "

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"        # force text so "1.200" is not parsed as a number
$cell.Value = "1.200"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E50").Value = "  -3.48%  "

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"        # force text so "1.146" is not parsed as a number
$cell.Value = "1.146"
$cell.Style = "Normal"          # restore default style (only the value changed)
$ws.Range("E51").Value = "  -0.93%  "
